$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825

$ws.Range("H96").Value = 630.7
$ws.Range("I96").Value = 301.45456
$ws.Range("J96").Value = 1033.1111
$ws.Range("K96").Value = 904.36368
$ws.Range("L96").Value = 3099.3333
$ws.Range("M96").Value = 468.63632
$ws.Range("N96").Value = -5845.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9450.924999999999
$ws.Range("I2").Value = 11629.219
$ws.Range("J2").Value = 737.75
$ws.Range("K2").Value = 11629.219
$ws.Range("L2").Value = 737.75
$ws.Range("M2").Value = -11516.219
$ws.Range("N2").Value = -963.75

$ws.Range("H27").Value = 4666.6665
$ws.Range("J27").Value = 4666.6665
$ws.Range("L27").Value = 4666.6665
$ws.Range("N27").Value = -5034.6665

$ws.Range("H32").Value = 2310471.2
$ws.Range("I32").Value = 2907133.5
$ws.Range("K32").Value = 2907133.5
$ws.Range("M32").Value = -2906846.5

$ws.Range("H110").Value = 1918.5385
$ws.Range("I110").Value = 1161.381
$ws.Range("K110").Value = 1161.381
$ws.Range("M110").Value = 883.6189999999999

$ws.Range("H116").Value = 9450.924999999999
$ws.Range("I116").Value = 11629.219
$ws.Range("J116").Value = 737.75
$ws.Range("K116").Value = 11629.219
$ws.Range("L116").Value = 737.75
$ws.Range("M116").Value = -9335.218999999999
$ws.Range("N116").Value = -5325.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9450.924999999999
$ws.Range("I3").Value = 11629.219
$ws.Range("J3").Value = 737.75
$ws.Range("K3").Value = 11629.219
$ws.Range("L3").Value = 737.75
$ws.Range("M3").Value = -11515.219
$ws.Range("N3").Value = -965.75

$ws.Range("H64").Value = 41666896
$ws.Range("I64").Value = 353
$ws.Range("J64").Value = 83333440
$ws.Range("K64").Value = 353
$ws.Range("L64").Value = 83333440
$ws.Range("M64").Value = -128
$ws.Range("N64").Value = -83333890

$ws.Range("H67").Value = 41666896
$ws.Range("I67").Value = 353
$ws.Range("J67").Value = 83333440
$ws.Range("K67").Value = 353
$ws.Range("L67").Value = 83333440
$ws.Range("M67").Value = 427
$ws.Range("N67").Value = -83335000

$ws.Range("H94").Value = 1229.8064
$ws.Range("I94").Value = 989.1667
$ws.Range("K94").Value = 989.1667
$ws.Range("M94").Value = -538.1667

$ws.Range("H105").Value = 1960
$ws.Range("J105").Value = 1900
$ws.Range("L105").Value = 1900
$ws.Range("N105").Value = -5394

$ws.Range("H107").Value = 987.625
$ws.Range("I107").Value = 1003.3333
$ws.Range("J107").Value = 940.5
$ws.Range("K107").Value = 1003.3333
$ws.Range("L107").Value = 940.5
$ws.Range("M107").Value = 916.6667
$ws.Range("N107").Value = -4780.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1341790.9
$ws.Range("I58").Value = 5732.905
$ws.Range("J58").Value = 3500038.2
$ws.Range("K58").Value = 5732.905
$ws.Range("L58").Value = 3500038.2
$ws.Range("M58").Value = -5529.905
$ws.Range("N58").Value = -3500444.2

$ws.Range("H136").Value = 1341790.9
$ws.Range("I136").Value = 5732.905
$ws.Range("J136").Value = 3500038.2
$ws.Range("K136").Value = 17198.715
$ws.Range("L136").Value = 10500114.6
$ws.Range("M136").Value = -14648.715
$ws.Range("N136").Value = -10505214.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1497.262
$ws.Range("I113").Value = 1202.6316
$ws.Range("J113").Value = 1740.6522
$ws.Range("K113").Value = 3607.8948
$ws.Range("L113").Value = 5221.9566
$ws.Range("M113").Value = -1437.8948
$ws.Range("N113").Value = -9561.9566

$ws.Range("H131").Value = 843.42
$ws.Range("I131").Value = 286.92307
$ws.Range("J131").Value = 926.5747
$ws.Range("K131").Value = 860.7692099999999
$ws.Range("L131").Value = 2779.7241
$ws.Range("M131").Value = 4179.23079
$ws.Range("N131").Value = -12859.7241

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 19816.4
$ws.Range("I113").Value = 2329.125
$ws.Range("J113").Value = 39801.855
$ws.Range("K113").Value = 2329.125
$ws.Range("L113").Value = 39801.855
$ws.Range("M113").Value = -159.125
$ws.Range("N113").Value = -44141.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 83334136
$ws.Range("I46").Value = 999.6667
$ws.Range("K46").Value = 999.6667
$ws.Range("M46").Value = -811.6667

$ws.Range("H61").Value = 3028.2856
$ws.Range("I61").Value = 2039.4
$ws.Range("J61").Value = 3927.2727
$ws.Range("K61").Value = 2039.4
$ws.Range("L61").Value = 3927.2727
$ws.Range("M61").Value = -1837.4
$ws.Range("N61").Value = -4331.2727

$ws.Range("H68").Value = 2599.7058
$ws.Range("I68").Value = 1294.5714
$ws.Range("J68").Value = 3513.3
$ws.Range("K68").Value = 1294.5714
$ws.Range("L68").Value = 3513.3
$ws.Range("M68").Value = -545.5714
$ws.Range("N68").Value = -5011.3

$ws.Range("H71").Value = 2599.7058
$ws.Range("I71").Value = 1294.5714
$ws.Range("J71").Value = 3513.3
$ws.Range("K71").Value = 6472.857
$ws.Range("L71").Value = 17566.5
$ws.Range("M71").Value = -2728.857
$ws.Range("N71").Value = -25054.5

$ws.Range("H82").Value = 3765.4546
$ws.Range("I82").Value = 1048.8889
$ws.Range("J82").Value = 5646.154
$ws.Range("K82").Value = 1048.8889
$ws.Range("L82").Value = 5646.154
$ws.Range("M82").Value = -687.8888999999999
$ws.Range("N82").Value = -6368.154

$ws.Range("H85").Value = 3765.4546
$ws.Range("I85").Value = 1048.8889
$ws.Range("J85").Value = 5646.154
$ws.Range("K85").Value = 1048.8889
$ws.Range("L85").Value = 5646.154
$ws.Range("M85").Value = 199.1111000000001
$ws.Range("N85").Value = -8142.154

$ws.Range("H113").Value = 3028.2856
$ws.Range("I113").Value = 2039.4
$ws.Range("J113").Value = 3927.2727
$ws.Range("K113").Value = 2039.4
$ws.Range("L113").Value = 3927.2727
$ws.Range("M113").Value = 130.5999999999999
$ws.Range("N113").Value = -8267.2727

$ws.Range("H132").Value = 11915061
$ws.Range("I132").Value = 17870466
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 53611398
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -53608868
$ws.Range("N132").Value = -17809.25

$ws.Range("H141").Value = 67928.625
$ws.Range("J141").Value = 67928.625
$ws.Range("L141").Value = 67928.625
$ws.Range("N141").Value = -78288.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 885104.2
$ws.Range("I132").Value = 4299.8335
$ws.Range("J132").Value = 1589747.6
$ws.Range("K132").Value = 12899.5005
$ws.Range("L132").Value = 4769242.800000001
$ws.Range("M132").Value = -10369.5005
$ws.Range("N132").Value = -4774302.800000001
